$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "KPIs" - overall KPI summary (row 2)
# ---------------------------------------------------------------------------
$wsKPIs = $wb.Worksheets.Item("KPIs")
$wsKPIs.Range("A2").Value = 7655.499999999999
$wsKPIs.Range("E2").Value = 159.4895833333333
$wsKPIs.Range("F2").Value = 0.631578947368421

# ---------------------------------------------------------------------------
# Sheet "Ventes Mensuelles" - monthly sales (rows 2-7, 2006 Jan..Jun)
# ---------------------------------------------------------------------------
$wsVentes = $wb.Worksheets.Item("Ventes Mensuelles")
$wsVentes.Range("C2").Value = 507.5
$wsVentes.Range("E2").Value = 34

$wsVentes.Range("C3").Value = 413.7
$wsVentes.Range("E3").Value = 31

$wsVentes.Range("C4").Value = 1331.2
$wsVentes.Range("E4").Value = 88

$wsVentes.Range("C5").Value = 2796.1
$wsVentes.Range("E5").Value = 192

$wsVentes.Range("C6").Value = 1027.6
$wsVentes.Range("E6").Value = 67

$wsVentes.Range("C7").Value = 1579.4
$wsVentes.Range("E7").Value = 115

# ---------------------------------------------------------------------------
# Sheet "Par Catégorie" - sales by category (rows 2-3)
# ---------------------------------------------------------------------------
$wsCategorie = $wb.Worksheets.Item("Par Catégorie")
$wsCategorie.Range("B2").Value = 4086
$wsCategorie.Range("D2").Value = 245

$wsCategorie.Range("B3").Value = 3569.5
$wsCategorie.Range("C3").Value = 32
$wsCategorie.Range("D3").Value = 282

# ---------------------------------------------------------------------------
# Sheet "Top Produits" - top products (rows 2-4)
# ---------------------------------------------------------------------------
$wsProduits = $wb.Worksheets.Item("Top Produits")
$wsProduits.Range("B2").Value = 4086
$wsProduits.Range("C2").Value = 245

$wsProduits.Range("B3").Value = 1787.5
$wsProduits.Range("C3").Value = 192
$wsProduits.Range("D3").Value = 32

$wsProduits.Range("B4").Value = 1782
$wsProduits.Range("C4").Value = 90

# ---------------------------------------------------------------------------
# Sheet "Par Pays" - sales by country (row 2)
# ---------------------------------------------------------------------------
$wsPays = $wb.Worksheets.Item("Par Pays")
$wsPays.Range("B2").Value = 7655.5

# ---------------------------------------------------------------------------
# Sheet "Employés" - sales by employee, re-ranked by TotalSales descending
# (rows 2-9)
# ---------------------------------------------------------------------------
$wsEmployes = $wb.Worksheets.Item("Employés")

$wsEmployes.Range("A2").Value = "Anne Hellung-Larsen"
$wsEmployes.Range("B2").Value = 2263.2
$wsEmployes.Range("C2").Value = 10
$wsEmployes.Range("D2").Value = 5

$wsEmployes.Range("A3").Value = "Mariya Sergienko"
$wsEmployes.Range("B3").Value = 1432.1
$wsEmployes.Range("C3").Value = 8
$wsEmployes.Range("D3").Value = 4

$wsEmployes.Range("A4").Value = "Nancy Freehafer"
$wsEmployes.Range("B4").Value = 1365
$wsEmployes.Range("C4").Value = 12
$wsEmployes.Range("D4").Value = 6

$wsEmployes.Range("A5").Value = "Jan Kotas"
$wsEmployes.Range("B5").Value = 776.3000000000001
$wsEmployes.Range("C5").Value = 6
$wsEmployes.Range("D5").Value = 3

$wsEmployes.Range("A6").Value = "Andrew Cencini"
$wsEmployes.Range("B6").Value = 699
$wsEmployes.Range("C6").Value = 4
$wsEmployes.Range("D6").Value = 3

$wsEmployes.Range("A7").Value = "Laura Giussani"
$wsEmployes.Range("B7").Value = 488.7
$wsEmployes.Range("C7").Value = 2
$wsEmployes.Range("D7").Value = 1

$wsEmployes.Range("A8").Value = "Robert Zare"
$wsEmployes.Range("B8").Value = 384.3
$wsEmployes.Range("C8").Value = 2
$wsEmployes.Range("D8").Value = 1

$wsEmployes.Range("A9").Value = "Michael Neipper"
$wsEmployes.Range("B9").Value = 246.9
$wsEmployes.Range("C9").Value = 4
$wsEmployes.Range("D9").Value = 2
